$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Osszeadas (addition) - fill in the Pass/Fail + message results
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F3").Value = "Fail"
$ws1.Range("G3").Value = "Hibás"

$ws1.Range("F4").Value = "Pass"

$ws1.Range("F5").Value = "Fail"
$ws1.Range("G5").Value = "Hibás"

$ws1.Range("F6").Value = "Pass"

$ws1.Activate()
$ws1.Range("I6").Select()

# ---------------------------------------------------------------------------
# Sheet 2: Kivonas (subtraction) - add the missing blank message cells
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("G3").Value = " "
$ws2.Range("G5").Value = " "

$ws2.Activate()
$ws2.Range("I7").Select()

# ---------------------------------------------------------------------------
# Sheet 3: Szorzas (multiplication) - populate the whole test table
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("B2").Value = "Bement 1 "
$ws3.Range("C2").Value = "Bemenet 2"
$ws3.Range("D2").Value = "Bemenet 3"
$ws3.Range("F2").Value = "Pass/Fail"
$ws3.Range("G2").Value = "Üzenet"

$ws3.Range("B3").Value = 10
$ws3.Range("C3").Value = 5
$ws3.Range("D3").Value = 50
$ws3.Range("F3").Value = "Pass"
$ws3.Range("G3").Value = " "

$ws3.Range("B4").Value = 20
$ws3.Range("C4").Value = 30
$ws3.Range("D4").Value = 600
$ws3.Range("F4").Value = "Pass"
$ws3.Range("G4").Value = " "

$ws3.Range("B5").Value = 70
$ws3.Range("C5").Value = 30
$ws3.Range("D5").Value = 2101
$ws3.Range("F5").Value = "Fail"
$ws3.Range("G5").Value = "Hibás"

$ws3.Range("B6").Value = 90
$ws3.Range("C6").Value = 60
$ws3.Range("D6").Value = 5400
$ws3.Range("F6").Value = "Pass"
$ws3.Range("G6").Value = " "

$ws3.Activate()
$ws3.Range("F3:G6").Select()

# ---------------------------------------------------------------------------
# Sheet 4: Osztas (division) - populate the whole test table
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("B2").Value = "Bement 1 "
$ws4.Range("C2").Value = "Bemenet 2"
$ws4.Range("D2").Value = "Bemenet 3"
$ws4.Range("F2").Value = "Pass/Fail"
$ws4.Range("G2").Value = "Üzenet"

$ws4.Range("B3").Value = 10
$ws4.Range("C3").Value = 5
$ws4.Range("D3").Value = 2
$ws4.Range("F3").Value = "Pass"
$ws4.Range("G3").Value = " "

$ws4.Range("B4").Value = 20
$ws4.Range("C4").Value = 30
$ws4.Range("D4").Value = 1
$ws4.Range("F4").Value = "Fail"
$ws4.Range("G4").Value = "Hibás"

$ws4.Range("B5").Value = 30
$ws4.Range("C5").Value = 30
$ws4.Range("D5").Value = 1
$ws4.Range("F5").Value = "Pass"
$ws4.Range("G5").Value = " "

$ws4.Range("B6").Value = 90
$ws4.Range("C6").Value = 60
$ws4.Range("D6").Value = 1.5
$ws4.Range("F6").Value = "Pass"
$ws4.Range("G6").Value = " "

$ws4.Activate()
$ws4.Range("L9").Select()
